$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOE")

# Change "Variable Type" values from "Continuous" to "Discrete"
$ws.Range("H2").Value = "Discrete"
$ws.Range("H3").Value = "Discrete"
$ws.Range("H4").Value = "Discrete"

# Clear out the old Mean/Standard Deviation columns (B, C) for rows 2-4
$ws.Range("B2:C4").ClearContents()

# Populate the new Max / Min / Step columns (D, E, F) for rows 2-4
$ws.Range("D2").Value = 130
$ws.Range("E2").Value = 70
$ws.Range("F2").Value = 0.06

$ws.Range("D3").Value = 0.13
$ws.Range("E3").Value = 0.07000000000000001
$ws.Range("F3").Value = 0.00005999999999999999

$ws.Range("D4").Value = 0.0013
$ws.Range("E4").Value = 0.0007
$ws.Range("F4").Value = 0.0000006
